$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-09-01T16:07:53"
$ws.Range("T4").Value = 71.63
$ws.Range("U4").Value = 80
$ws.Range("V4").Value = 50.07
$ws.Range("W4").Value = 39.86
$ws.Range("Y4").Value = 34.26
$ws.Range("Z4").Value = 34.63
$ws.Range("T6").Value = -3.15
$ws.Range("U6").Value = -3.36
$ws.Range("V6").Value = -2.05
$ws.Range("W6").Value = -1.63
$ws.Range("Y6").Value = -1.34
$ws.Range("Z6").Value = -0.93
$ws.Range("T9").Value = 70.95
$ws.Range("U9").Value = 80.23
$ws.Range("V9").Value = 49.92
$ws.Range("W9").Value = 39.79
$ws.Range("Z9").Value = 36.55
$ws.Range("T11").Value = -3.83
$ws.Range("U11").Value = -3.13
$ws.Range("V11").Value = -2.2
$ws.Range("W11").Value = -1.71
$ws.Range("Z11").Value = 0.99
$ws.Range("T14").Value = 70.95
$ws.Range("U14").Value = 80.23
$ws.Range("V14").Value = 49.92
$ws.Range("W14").Value = 45.23
$ws.Range("X14").Value = 180
$ws.Range("Y14").Value = 219.59
$ws.Range("Z14").Value = 116.95
$ws.Range("W15").Value = 5.44
$ws.Range("X15").Value = 144.68
$ws.Range("Y15").Value = 185.2
$ws.Range("Z15").Value = 80.37
$ws.Range("T16").Value = -3.83
$ws.Range("U16").Value = -3.13
$ws.Range("V16").Value = -2.2
$ws.Range("W16").Value = -1.71
$ws.Range("Z16").Value = 1.02
$ws.Range("T19").Value = 71.15000000000001
$ws.Range("U19").Value = 79.62
$ws.Range("V19").Value = 49.73
$ws.Range("W19").Value = 39.63
$ws.Range("T21").Value = -3.63
$ws.Range("U21").Value = -3.74
$ws.Range("V21").Value = -2.39
$ws.Range("W21").Value = -1.86
$ws.Range("T24").Value = 71.15000000000001
$ws.Range("U24").Value = 79.62
$ws.Range("V24").Value = 49.73
$ws.Range("W24").Value = 39.63
$ws.Range("T26").Value = -3.63
$ws.Range("U26").Value = -3.74
$ws.Range("V26").Value = -2.39
$ws.Range("W26").Value = -1.86
$ws.Range("T29").Value = 70.68000000000001
$ws.Range("U29").Value = 79.23999999999999
$ws.Range("V29").Value = 49.4
$ws.Range("W29").Value = 39.41
$ws.Range("X29").Value = 34.77
$ws.Range("Z29").Value = 34.83
$ws.Range("T31").Value = -4.1
$ws.Range("U31").Value = -4.12
$ws.Range("V31").Value = -2.72
$ws.Range("W31").Value = -2.09
$ws.Range("X31").Value = -1.46
$ws.Range("Z31").Value = -0.73
$ws.Range("T34").Value = 70.88
$ws.Range("U34").Value = 80.62
$ws.Range("V34").Value = 50.21
$ws.Range("W34").Value = 45.46
$ws.Range("X34").Value = 180.21
$ws.Range("Y34").Value = 219.99
$ws.Range("Z34").Value = 117.8
$ws.Range("W35").Value = 5.44
$ws.Range("X35").Value = 144.68
$ws.Range("Y35").Value = 185.2
$ws.Range("Z35").Value = 80.37
$ws.Range("T36").Value = -3.9
$ws.Range("U36").Value = -2.74
$ws.Range("V36").Value = -1.91
$ws.Range("W36").Value = -1.48
$ws.Range("Z36").Value = 1.87
$ws.Range("T39").Value = 71.63
$ws.Range("U39").Value = 80
$ws.Range("V39").Value = 50.07
$ws.Range("W39").Value = 39.86
$ws.Range("Y39").Value = 34.26
$ws.Range("Z39").Value = 34.63
$ws.Range("T41").Value = -3.15
$ws.Range("U41").Value = -3.36
$ws.Range("V41").Value = -2.05
$ws.Range("W41").Value = -1.63
$ws.Range("Y41").Value = -1.34
$ws.Range("Z41").Value = -0.93
$ws.Range("T44").Value = 75.68000000000001
$ws.Range("U44").Value = 84.08
$ws.Range("V44").Value = 52.61
$ws.Range("W44").Value = 42
$ws.Range("X44").Value = 36.53
$ws.Range("T46").Value = 0.9
$ws.Range("U46").Value = 0.72
$ws.Range("V46").Value = 0.49
$ws.Range("W46").Value = 0.51
$ws.Range("X46").Value = 0.3
$ws.Range("T49").Value = 74.33
$ws.Range("U49").Value = 86.73999999999999
$ws.Range("V49").Value = 54.01
$ws.Range("W49").Value = 42.91
$ws.Range("T51").Value = -0.45
$ws.Range("U51").Value = 3.38
$ws.Range("V51").Value = 1.89
$ws.Range("W51").Value = 1.42
$ws.Range("T54").Value = 75.61
$ws.Range("U54").Value = 84.45999999999999
$ws.Range("V54").Value = 52.12
$ws.Range("W54").Value = 41.75
$ws.Range("X54").Value = 37.13
$ws.Range("Y54").Value = 36.51
$ws.Range("T56").Value = 0.83
$ws.Range("U56").Value = 1.1
$ws.Range("V56").Value = 0
$ws.Range("W56").Value = 0.25
$ws.Range("X56").Value = 0.89
$ws.Range("Y56").Value = 0.91
$ws.Range("T59").Value = 78.72
$ws.Range("U59").Value = 87.09999999999999
$ws.Range("V59").Value = 54.46
$ws.Range("W59").Value = 43.59
$ws.Range("T61").Value = 3.94
$ws.Range("U61").Value = 3.75
$ws.Range("V61").Value = 2.34
$ws.Range("W61").Value = 2.09
$ws.Range("T64").Value = 79.72
$ws.Range("U64").Value = 88.12
$ws.Range("V64").Value = 55.03
$ws.Range("W64").Value = 44.05
$ws.Range("T66").Value = 4.94
$ws.Range("U66").Value = 4.76
$ws.Range("V66").Value = 2.92
$ws.Range("W66").Value = 2.55
$ws.Range("T69").Value = 80.48999999999999
$ws.Range("U69").Value = 88.68000000000001
$ws.Range("V69").Value = 55.27
$ws.Range("W69").Value = 44.52
$ws.Range("T71").Value = 5.72
$ws.Range("U71").Value = 5.32
$ws.Range("V71").Value = 3.15
$ws.Range("W71").Value = 3.03
$ws.Range("T74").Value = 78.06
$ws.Range("U74").Value = 86.2
$ws.Range("V74").Value = 54.01
$ws.Range("W74").Value = 43.23
$ws.Range("T76").Value = 3.28
$ws.Range("U76").Value = 2.84
$ws.Range("V76").Value = 1.89
$ws.Range("W76").Value = 1.73
$ws.Range("T79").Value = 78.48
$ws.Range("U79").Value = 86.59999999999999
$ws.Range("V79").Value = 54.28
$ws.Range("W79").Value = 43.47
$ws.Range("Z79").Value = 36.88
$ws.Range("T81").Value = 3.7
$ws.Range("U81").Value = 3.24
$ws.Range("V81").Value = 2.17
$ws.Range("W81").Value = 1.97
$ws.Range("T84").Value = 75.61
$ws.Range("U84").Value = 84.89
$ws.Range("V84").Value = 51.15
$ws.Range("W84").Value = 41
$ws.Range("T86").Value = 0.83
$ws.Range("U86").Value = 1.53
$ws.Range("V86").Value = -0.97
$ws.Range("W86").Value = -0.49
$ws.Range("T89").Value = 70.68000000000001
$ws.Range("U89").Value = 79.23999999999999
$ws.Range("V89").Value = 49.4
$ws.Range("W89").Value = 39.41
$ws.Range("X89").Value = 34.77
$ws.Range("Z89").Value = 34.83
$ws.Range("T91").Value = -4.1
$ws.Range("U91").Value = -4.12
$ws.Range("V91").Value = -2.72
$ws.Range("W91").Value = -2.09
$ws.Range("X91").Value = -1.46
$ws.Range("Z91").Value = -0.73
